$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$newRows = @(
    @{ Row = 245; Date = 44319; B = 0; C = 1; D = 62.34413965087282 },
    @{ Row = 246; Date = 44320; B = 0; C = 1; D = 62.34413965087282 },
    @{ Row = 247; Date = 44321; B = 0; C = 1; D = 62.34413965087282 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy the formatting of the last existing data row (244) into column A of the new row
    $ws.Range("A244").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
